# Update time_taken timestamps (column F) on the "data" sheet to reflect
# the re-run query time (13:40:54.xxx -> 14:21:28.xxx)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:21:28.052056"
$ws1.Range("F3").Value = "2021-10-05 14:21:28.052064"
$ws1.Range("F4").Value = "2021-10-05 14:21:28.052067"
$ws1.Range("F5").Value = "2021-10-05 14:21:28.052069"
$ws1.Range("F6").Value = "2021-10-05 14:21:28.052072"
$ws1.Range("F7").Value = "2021-10-05 14:21:28.052075"
$ws1.Range("F8").Value = "2021-10-05 14:21:28.052077"
$ws1.Range("F9").Value = "2021-10-05 14:21:28.052079"
$ws1.Range("F10").Value = "2021-10-05 14:21:28.052082"
$ws1.Range("F11").Value = "2021-10-05 14:21:28.052084"
$ws1.Range("F12").Value = "2021-10-05 14:21:28.052087"
$ws1.Range("F13").Value = "2021-10-05 14:21:28.052089"
$ws1.Range("F14").Value = "2021-10-05 14:21:28.052092"
$ws1.Range("F15").Value = "2021-10-05 14:21:28.052094"
$ws1.Range("F16").Value = "2021-10-05 14:21:28.052096"
$ws1.Range("F17").Value = "2021-10-05 14:21:28.052099"
$ws1.Range("F18").Value = "2021-10-05 14:21:28.052102"
$ws1.Range("F19").Value = "2021-10-05 14:21:28.052104"
$ws1.Range("F20").Value = "2021-10-05 14:21:28.052106"
$ws1.Range("F21").Value = "2021-10-05 14:21:28.052109"
$ws1.Range("F22").Value = "2021-10-05 14:21:28.052111"
$ws1.Range("F23").Value = "2021-10-05 14:21:28.052114"
$ws1.Range("F24").Value = "2021-10-05 14:21:28.052116"
$ws1.Range("F25").Value = "2021-10-05 14:21:28.052119"
$ws1.Range("F26").Value = "2021-10-05 14:21:28.052122"
$ws1.Range("F27").Value = "2021-10-05 14:21:28.052124"
$ws1.Range("F28").Value = "2021-10-05 14:21:28.052127"
$ws1.Range("F29").Value = "2021-10-05 14:21:28.052129"
$ws1.Range("F30").Value = "2021-10-05 14:21:28.052131"
$ws1.Range("F31").Value = "2021-10-05 14:21:28.052134"
$ws1.Range("F32").Value = "2021-10-05 14:21:28.052136"
$ws1.Range("F33").Value = "2021-10-05 14:21:28.052139"
$ws1.Range("F34").Value = "2021-10-05 14:21:28.052142"
$ws1.Range("F35").Value = "2021-10-05 14:21:28.052144"
$ws1.Range("F36").Value = "2021-10-05 14:21:28.052146"
$ws1.Range("F37").Value = "2021-10-05 14:21:28.052149"
$ws1.Range("F38").Value = "2021-10-05 14:21:28.052151"
$ws1.Range("F39").Value = "2021-10-05 14:21:28.052154"
$ws1.Range("F40").Value = "2021-10-05 14:21:28.052156"
$ws1.Range("F41").Value = "2021-10-05 14:21:28.052158"
$ws1.Range("F42").Value = "2021-10-05 14:21:28.052161"
$ws1.Range("F43").Value = "2021-10-05 14:21:28.052164"
$ws1.Range("F44").Value = "2021-10-05 14:21:28.052166"
$ws1.Range("F45").Value = "2021-10-05 14:21:28.052169"
$ws1.Range("F46").Value = "2021-10-05 14:21:28.052171"
$ws1.Range("F47").Value = "2021-10-05 14:21:28.052174"
$ws1.Range("F48").Value = "2021-10-05 14:21:28.052176"
$ws1.Range("F49").Value = "2021-10-05 14:21:28.052178"
$ws1.Range("F50").Value = "2021-10-05 14:21:28.052181"
$ws1.Range("F51").Value = "2021-10-05 14:21:28.052183"
$ws1.Range("F52").Value = "2021-10-05 14:21:28.052186"
$ws1.Range("F53").Value = "2021-10-05 14:21:28.052188"
$ws1.Range("F54").Value = "2021-10-05 14:21:28.052191"
$ws1.Range("F55").Value = "2021-10-05 14:21:28.052193"
$ws1.Range("F56").Value = "2021-10-05 14:21:28.052196"

# Add a new "metadata" worksheet after "data", describing the panel query
# that produced the data sheet.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "metadata"

# Copy the header style (bold, centered, bordered) used in "data" B1 to the
# header row of "metadata", and the index-column style used in "data" A2 to
# "metadata" A2.
$ws1.Range("B1").Copy()
foreach ($col in @("B","C","D","E","F","G")) {
    $newSheet.Range($col + "1").PasteSpecial(-4122)
}
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Lysosomal storage disorder"
$newSheet.Range("C2").Value = 529
# data_version "1.74" must stay a text value (matches "data" sheet's style
# of storing these identifiers as text), not be coerced into a number.
$newSheet.Range("D2").Value = "'1.74"
$newSheet.Range("D2").ClearFormats()
$newSheet.Range("E2").Value = "2021-06-14T11:39:49.478955Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:28.048930"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/529/?format=json"
